# Weekly data refresh: insert this week's new price record for
# Ciruela / Angeleno / Primera (O'Higgins) ahead of the existing
# historical rows, pushing the old rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 233; everything from the old row 233
# downward (through the old last row 242) shifts down to 234..243.
$ws.Rows.Item(233).Insert()

# Populate the newly inserted row 233 with this week's record.
$ws.Cells.Item(233, 1).Value2 = 10
$ws.Cells.Item(233, 2).Value2 = "Vega Modelo de Temuco"
$ws.Cells.Item(233, 3).Value2 = "La Araucanía"
$ws.Cells.Item(233, 4).Value2 = 44753
$ws.Cells.Item(233, 5).Value2 = 9
$ws.Cells.Item(233, 6).Value2 = "Fruta"
$ws.Cells.Item(233, 7).Value2 = 100103
$ws.Cells.Item(233, 8).Value2 = "Frutos de hueso (carozo)"
$ws.Cells.Item(233, 9).Value2 = 100103002
$ws.Cells.Item(233, 10).Value2 = "Ciruela"
$ws.Cells.Item(233, 11).Value2 = "Angeleno"
$ws.Cells.Item(233, 12).Value2 = "Primera"
$ws.Cells.Item(233, 13).Value2 = 80
$ws.Cells.Item(233, 14).Value2 = 10000
$ws.Cells.Item(233, 15).Value2 = 10000
$ws.Cells.Item(233, 16).Value2 = 10000
$ws.Cells.Item(233, 17).Value2 = "`$/bandeja 18 kilos granel"
$ws.Cells.Item(233, 18).Value2 = "Región de O'Higgins"
$ws.Cells.Item(233, 19).Value2 = 556
$ws.Cells.Item(233, 20).Value2 = 18
